{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is `async (context) => { ... }`.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that exactly reads \"Data, Technology and Strategy Consulting\"\n// (the line directly under the Siege Analytics/PARTNER heading).\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.trim() === \"Data, Technology and Strategy Consulting\") {\n    anchor = p;\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the 'Data, Technology and Strategy Consulting' paragraph.\");\n}\n\n// Insert the three new bullet paragraphs right after it, preserving order.\nconst newLines = [\n  \"\\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n  \"\\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n  \"\\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n];\n\nlet insertAfter = anchor;\nfor (const line of newLines) {\n  insertAfter = insertAfter.insertParagraph(line, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that exactly reads \"Data, Technology and Strategy Consulting\"\n# (the line directly under the Siege Analytics/PARTNER heading).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq 'Data, Technology and Strategy Consulting') {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Data, Technology and Strategy Consulting' paragraph.\"\n}\n\n# New bullet lines to insert directly after it, in order.\n$lines = @(\n    [char]0x2022 + ' Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters',\n    [char]0x2022 + ' Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States',\n    [char]0x2022 + ' Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis'\n)\n\n$cur = $target\nforeach ($line in $lines) {\n    $cur.Range.InsertParagraphAfter()\n    $cur = $cur.Next()\n    $cur.Range.Text = $line\n}\n"}
